$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1099
$ws.Range("I97").Value = 450
$ws.Range("J97").Value = 2397
$ws.Range("K97").Value = 1350
$ws.Range("L97").Value = 7191
$ws.Range("M97").Value = -854
$ws.Range("N97").Value = -8183
$ws.Range("H106").Value = 2360.9473
$ws.Range("I106").Value = 1462.2941
$ws.Range("K106").Value = 1462.2941
$ws.Range("M106").Value = -831.2941000000001
$ws.Range("H107").Value = 1739.92
$ws.Range("I107").Value = 1516.375
$ws.Range("J107").Value = 2137.3333
$ws.Range("K107").Value = 1516.375
$ws.Range("L107").Value = 2137.3333
$ws.Range("M107").Value = 403.625
$ws.Range("N107").Value = -5977.3333
$ws.Range("H111").Value = 4355.1113
$ws.Range("I111").Value = 2413
$ws.Range("J111").Value = 8239.333000000001
$ws.Range("K111").Value = 7239
$ws.Range("L111").Value = 24717.999
$ws.Range("M111").Value = -4172
$ws.Range("N111").Value = -30851.999
$ws.Range("H137").Value = 12600.695
$ws.Range("J137").Value = 21360.455
$ws.Range("L137").Value = 64081.36500000001
$ws.Range("N137").Value = -69181.36500000001
$ws.Range("H138").Value = 6490.069
$ws.Range("I138").Value = 2555.5
$ws.Range("J138").Value = 6781.5186
$ws.Range("K138").Value = 7666.5
$ws.Range("L138").Value = 20344.5558
$ws.Range("M138").Value = -2526.5
$ws.Range("N138").Value = -30624.5558
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 39999
$ws.Range("J61").Value = 39999
$ws.Range("L61").Value = 39999
$ws.Range("N61").Value = -40625
$ws.Range("H132").Value = 89709.14
$ws.Range("J132").Value = 89709.14
$ws.Range("L132").Value = 89709.14
$ws.Range("N132").Value = -99829.14

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 28055.5
$ws.Range("J43").Value = 28055.5
$ws.Range("L43").Value = 28055.5
$ws.Range("N43").Value = -28423.5
$ws.Range("H50").Value = 34999
$ws.Range("J50").Value = 34999
$ws.Range("L50").Value = 34999
$ws.Range("N50").Value = -36249
$ws.Range("H51").Value = 29999
$ws.Range("J51").Value = 29999
$ws.Range("L51").Value = 29999
$ws.Range("N51").Value = -31471
$ws.Range("H59").Value = 40057.5
$ws.Range("I59").Value = 40000
$ws.Range("J59").Value = 40115
$ws.Range("K59").Value = 40000
$ws.Range("L59").Value = 40115
$ws.Range("M59").Value = -38855
$ws.Range("N59").Value = -42405
$ws.Range("H60").Value = 45000
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H61").Value = 29999
$ws.Range("J61").Value = 29999
$ws.Range("L61").Value = 29999
$ws.Range("N61").Value = -30695
$ws.Range("H101").Value = 28055.5
$ws.Range("J101").Value = 28055.5
$ws.Range("L101").Value = 28055.5
$ws.Range("N101").Value = -34545.5
$ws.Range("H103").Value = 38799.4
$ws.Range("I103").Value = 26998.5
$ws.Range("J103").Value = 46666.668
$ws.Range("K103").Value = 26998.5
$ws.Range("L103").Value = 46666.668
$ws.Range("M103").Value = -25826.5
$ws.Range("N103").Value = -49010.668
$ws.Range("H111").Value = 77000
$ws.Range("J111").Value = 77000
$ws.Range("L111").Value = 77000
$ws.Range("N111").Value = -85180
$ws.Range("H133").Value = 59320.5
$ws.Range("J133").Value = 59320.5
$ws.Range("L133").Value = 59320.5
$ws.Range("N133").Value = -64380.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 149988.33
$ws.Range("J37").Value = 149988.33
$ws.Range("L37").Value = 449964.99
$ws.Range("N37").Value = -450188.99
$ws.Range("H113").Value = 1494.25
$ws.Range("J113").Value = 1494.25
$ws.Range("L113").Value = 4482.75
$ws.Range("N113").Value = -8822.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 620.2381
$ws.Range("J97").Value = 717.8570999999999
$ws.Range("L97").Value = 717.8570999999999
$ws.Range("N97").Value = -1709.8571
$ws.Range("H102").Value = 2524.75
$ws.Range("I102").Value = 1999.5
$ws.Range("K102").Value = 1999.5
$ws.Range("M102").Value = -377.5
$ws.Range("H133").Value = 81094.57000000001
$ws.Range("J133").Value = 81094.57000000001
$ws.Range("L133").Value = 81094.57000000001
$ws.Range("N133").Value = -91214.57000000001
$ws.Range("H139").Value = 70326
$ws.Range("J139").Value = 70326
$ws.Range("L139").Value = 70326
$ws.Range("N139").Value = -80606

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2161
$ws.Range("I46").Value = 979
$ws.Range("K46").Value = 979
$ws.Range("M46").Value = -791
$ws.Range("H61").Value = 2199
$ws.Range("I61").Value = 2199
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2199
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1997
$ws.Range("N61").ClearContents()
$ws.Range("H69").Value = 51273.5
$ws.Range("J69").Value = 49032.668
$ws.Range("L69").Value = 49032.668
$ws.Range("N69").Value = -50654.668
$ws.Range("H72").Value = 51273.5
$ws.Range("J72").Value = 49032.668
$ws.Range("L72").Value = 147098.004
$ws.Range("N72").Value = -155210.004
$ws.Range("H113").Value = 2199
$ws.Range("I113").Value = 2199
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2199
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -29
$ws.Range("N113").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4454.273
$ws.Range("I62").Value = 4666.6665
$ws.Range("K62").Value = 4666.6665
$ws.Range("M62").Value = -4042.6665
$ws.Range("H65").Value = 4454.273
$ws.Range("I65").Value = 4666.6665
$ws.Range("K65").Value = 23333.3325
$ws.Range("M65").Value = -20213.3325
$ws.Range("H81").Value = 11186.4375
$ws.Range("J81").Value = 16899.4
$ws.Range("L81").Value = 33798.8
$ws.Range("N81").Value = -35920.8
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H84").Value = 11186.4375
$ws.Range("J84").Value = 16899.4
$ws.Range("L84").Value = 168994
$ws.Range("N84").Value = -179602
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H113").Value = 910.2727
$ws.Range("I113").Value = 993.6667
$ws.Range("J113").Value = 731.5714
$ws.Range("K113").Value = 2981.0001
$ws.Range("L113").Value = 2194.7142
$ws.Range("M113").Value = -811.0001000000002
$ws.Range("N113").Value = -6534.7142
$ws.Range("H122").Value = 1900.3846
$ws.Range("I122").Value = 1922.2
$ws.Range("K122").Value = 5766.6
$ws.Range("M122").Value = -3316.6
